# atualizei dados da add
# Update faturamento_diario (Sheet1): correct two existing June totals and
# insert a missing daily record (day 10) for June/2025, which pushes every
# subsequent row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct previously wrong totals for day 6 and day 9 (June/2025).
$ws.Range("B6").Value = 12978.82
$ws.Range("B7").Value = 10572.4

# Insert a new row at position 8 for the missing day 10 (June/2025) record;
# this shifts all rows from 8..67 down to 9..68.
$ws.Rows.Item(8).Insert()

$ws.Range("A8").Value = 10
$ws.Range("B8").Value = 4076.75
$ws.Range("C8").Value = 6
$ws.Range("D8").Value = 2025
$ws.Range("E8").Value = "06/2025"
